$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.427.15"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "1.824.87"
$ws.Range("E3").Value = "  -0.32%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").Value = "'315.50"
$ws.Range("E5").Value = "  -1.00%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").Value = "'0.5126"
$ws.Range("E7").Value = "  -3.32%  "

$ws.Range("D8").Value = "'0.3918"
$ws.Range("E8").Value = "  -1.39%  "

$ws.Range("D9").Value = "'0.07647"
$ws.Range("E9").Value = "  +0.99%  "

$ws.Range("D10").Value = "'1.108"
$ws.Range("E10").Value = "  +0.33%  "

$ws.Range("D11").Value = "'41.64"
$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("D12").Value = "'21.04"
$ws.Range("E12").Value = "  +1.42%  "

$ws.Range("D13").Value = "'6.274"
$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("D14").Value = "'1.002"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").Value = "'7.519"
$ws.Range("E15").Value = "  -1.23%  "

$ws.Range("D16").Value = "1.824.77"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Value = "'93.48"
$ws.Range("E17").Value = "  +3.97%  "

$ws.Range("D18").Value = "'0.00001099"
$ws.Range("E18").Value = "  +2.55%  "

$ws.Range("D19").Value = "'0.06673"
$ws.Range("E19").Value = "  +1.21%  "

$ws.Range("D20").Value = "'17.64"
$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").Value = "'6.138"
$ws.Range("E22").Value = "  +1.19%  "

$ws.Range("D23").Value = "28.474.23"
$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").Value = "'11.15"
$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").Value = "'2.255"
$ws.Range("E25").Value = "  +7.94%  "

$ws.Range("D26").Value = "'20.69"
$ws.Range("E26").Value = "  +0.84%  "

$ws.Range("D27").Value = "'156.98"
$ws.Range("E27").Value = "  +0.16%  "

$ws.Range("D28").Value = "2.036.21"
$ws.Range("E28").Value = "  +0.20%  "

$ws.Range("D29").Value = "'2.375"
$ws.Range("E29").Value = "  -1.95%  "

$ws.Range("D30").Value = "'124.31"
$ws.Range("E30").Value = "  +0.67%  "

$ws.Range("D31").Value = "'1.108"
$ws.Range("E31").Value = "  -0.46%  "

$ws.Range("D32").Value = "'0.1087"
$ws.Range("E32").Value = "  -1.19%  "

$ws.Range("D33").Value = "'5.630"
$ws.Range("E33").Value = "  +0.35%  "

$ws.Range("D34").Value = "'3.655"
$ws.Range("E34").Value = "  -1.06%  "

$ws.Range("D35").Value = "'0.07026"
$ws.Range("E35").Value = "  -4.01%  "

$ws.Range("D36").Value = "'0.2202"
$ws.Range("E36").Value = "  -1.97%  "

$ws.Range("D37").Value = "'8.876"
$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").Value = "'0.02320"
$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("D39").Value = "'5.145"
$ws.Range("E39").Value = "  -1.64%  "

$ws.Range("D40").Value = "'0.6256"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").Value = "'11.19"
$ws.Range("E41").Value = "  -1.36%  "

$ws.Range("E42").Value = "  -1.90%  "

$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "'1.391"
$ws.Range("E44").Value = "  -1.72%  "

$ws.Range("D45").Value = "'13.45"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").Value = "'0.5887"
$ws.Range("E46").Value = "  +1.20%  "

$ws.Range("D47").Value = "'3.707"
$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("D48").Value = "'125.14"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("D49").Value = "'1.969"
$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").Value = "'1.195"
$ws.Range("E50").Value = "  +0.34%  "

$ws.Range("D51").Value = "'0.06926"
$ws.Range("E51").Value = "  +0.49%  "
